$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Correct/total marks updates
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
